$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'33.800.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.60%  '
$ws.Range("D3").Value = "'1.777.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.19%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'225.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.61%  '
$ws.Range("D6").Value = "'0.559"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.52%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'30.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").Value = "'46.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").Value = "'0.278"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.31%  '
$ws.Range("E11").Value = '  +3.57%  '
$ws.Range("D12").Value = "'0.0922"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").Value = "'2.030.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.05%  '
$ws.Range("D14").Value = "'1.774.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.99%  '
$ws.Range("D15").Value = "'0.625"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").Value = "'33.727.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.52%  '
$ws.Range("D17").Value = "'10.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.72%  '
$ws.Range("D19").Value = "'68.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.23%  '
$ws.Range("D20").Value = "'251.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = "'0.0₃0739"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = "'10.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("D25").Value = "'2.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").Value = "'159.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = "'16.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.81%  '
$ws.Range("E28").Value = '  +1.13%  '
$ws.Range("D29").Value = "'6.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").Value = "'0.0514"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("D33").Value = "'1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("D34").Value = "'3.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.16%  '
$ws.Range("D35").Value = "'1.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.18%  '
$ws.Range("D36").Value = "'1.480.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("D37").Value = "'1.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("D38").Value = "'0.631"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.40%  '
$ws.Range("E39").Value = '  +2.17%  '
$ws.Range("D40").Value = "'83.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("E41").Value = '  +2.11%  '
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("D43").Value = "'0.884"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.60%  '
$ws.Range("E44").Value = '  +2.01%  '
$ws.Range("D45").Value = "'0.0508"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.31%  '
$ws.Range("E46").Value = '  +3.50%  '
$ws.Range("D47").Value = "'1.928.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.62%  '
$ws.Range("D48").Value = "'5.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("D50").Value = "'11.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +14.84%  '
$ws.Range("D51").Value = "'50.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.43%  '
